$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the outdated row (Dia 2025-04-03 / SIND.DAS EMPRESAS... / produto 000782)
$ws.Rows(2).Delete()

# Renumber column A sequentially (1..6, then 0) to match the refreshed dataset
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 0
